# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 442 (shifting existing rows 442:587 down to 443:588),
# then populate the new row as a copy of the (original) row 442 data but with an updated date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 442; existing rows 442-587 shift down to 443-588.
$ws.Rows("442:442").Insert()

# Row 443 now contains what used to be row 442's data. Duplicate it into the new row 442.
$ws.Range("A443:T443").Copy()
$ws.Range("A442").PasteSpecial()
$ws.Application.CutCopyMode = $false

# Update the date of the newly inserted row to the new reading's date.
$ws.Range("D442").Value = 44985
